# Apply updated cryptocurrency data per diff (values, plus row 22/23 and 50/51 swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.920.82'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.80%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.633.65'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.98%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.47'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.90%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.517'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.11%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.54'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.21%  '

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.51%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +0.88%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.867.94'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.00%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.623.31'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.23%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.564'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.24%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.33'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +16.95%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.941.02'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.77%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.13%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.08'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.18%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '242.19'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.05%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.58%  '

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.17%  '

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.74%  '

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'Avalanche'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.81'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +4.18%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +3.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.70'

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.51'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.45%  '

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.78%  '

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.66%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.24%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0488'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.12%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +4.16%  '

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.38%  '

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.18%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.426.51'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.47%  '

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +4.87%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.63%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.62%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.30'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.17%  '

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.56%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.11'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +13.36%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.552'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.29%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.01'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.829'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.53%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.64%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.36%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '52.89'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.82%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.12%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.775.00'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +2.00%  '

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0114'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +10.90%  '

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Quant'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '90.91'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +4.98%  '
